$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Title paragraph: collapse the three runs into a single run with
#    new wording, re-centre it, switch off bold/underline, bump the
#    size to 22pt (sz 44) and switch the font to Liberation Serif.
#    We rebuild paragraph 1 via InsertXML (supplying a trailing
#    throw-away paragraph so the engine actually applies the new
#    <w:pPr>/<w:rPr> instead of keeping the old ones), then delete the
#    throw-away paragraph that InsertXML forces us to add.
# ---------------------------------------------------------------------

$titlePara = $d.Paragraphs(1).Range

$titleXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Liberation Serif" w:hAnsi="Liberation Serif" w:cs="Liberation Serif"/><w:sz w:val="44"/><w:szCs w:val="36"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Liberation Serif" w:hAnsi="Liberation Serif" w:cs="Liberation Serif"/><w:sz w:val="44"/><w:szCs w:val="36"/></w:rPr><w:t>Sketchs principaux et description</w:t></w:r></w:p><w:p><w:r><w:t>TMP_PLACEHOLDER</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$titlePara.InsertXML($titleXml)

# Remove the throw-away paragraph that was appended above.
$d.Paragraphs(2).Range.Delete()

# ---------------------------------------------------------------------
# 2) Re-create the "_GoBack" bookmark around the (now rewritten) title
#    paragraph. Word only ever keeps a single "_GoBack" bookmark, so
#    adding the new one transparently removes the old one that used to
#    sit around the second picture.
# ---------------------------------------------------------------------

$titleEnd = $d.Paragraphs(1).Range.End
$d.Bookmarks.Add("_GoBack", $d.Range(0, $titleEnd))

# ---------------------------------------------------------------------
# 3) Mark the second picture's run as <w:noProof/> (matches the other
#    picture runs already in the document).
# ---------------------------------------------------------------------

for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $shape = $d.InlineShapes($i)
    if ($shape.Width -eq 337.5 -and $shape.Height -eq 300) {
        $shape.Range.NoProofing = 1
    }
}
